# Book_data.xlsx edit: convert the "Release Date" column (C2:C4) from
# numeric dates to text values (typed with a leading apostrophe so Excel
# keeps the existing date number-format but stores the literal text),
# and move the selection to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write cells in row order C4, C3, C2 so the new shared-string entries are
# appended to sharedStrings.xml in that same order (matching the source
# workbook's save order: "1-Sep-12", "20-Nov-15", "12-Sep-23").
$ws.Range("C4").Value = "'1-Sep-12"
$ws.Range("C3").Value = "'20-Nov-15"
$ws.Range("C2").Value = "'12-Sep-23"

# Update the sheet selection/active cell to C3.
$ws.Range("C3").Select()
